# Update res_line/loading_percent values for case with 380 kV (Case_2_61)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new value
$updates = @(
    ,@(2, 2, 12.48608613848814)
    ,@(2, 3, 8.814538854320366)
    ,@(2, 5, 21.82773851148199)
    ,@(2, 6, 37.46172020787205)
    ,@(2, 7, 3.611554555625247)
    ,@(2, 9, 18.00890376160695)
    ,@(2, 10, 7.587886992288661)
    ,@(2, 13, 19.38967581523424)
    ,@(2, 15, 18.98200087479208)
    ,@(3, 2, 11.82612156877964)
    ,@(3, 3, 8.329607497461636)
    ,@(3, 5, 21.79157096143352)
    ,@(3, 6, 37.46203259756178)
    ,@(3, 7, 3.613523555236858)
    ,@(3, 9, 18.16207663472052)
    ,@(3, 10, 7.608300620032479)
    ,@(3, 13, 19.12769792378552)
    ,@(3, 15, 19.09758433813116)
    ,@(4, 2, 11.40123874880114)
    ,@(4, 3, 8.015638822159712)
    ,@(4, 5, 21.77337713718551)
    ,@(4, 6, 37.47361663186551)
    ,@(4, 7, 3.614795528525732)
    ,@(4, 9, 18.26117053013564)
    ,@(4, 10, 7.621493058583279)
    ,@(4, 13, 18.96732399825024)
    ,@(4, 15, 19.17469531113005)
    ,@(5, 2, 11.2233042509979)
    ,@(5, 3, 7.883682204362522)
    ,@(5, 5, 21.76697690964499)
    ,@(5, 6, 37.48119787061782)
    ,@(5, 7, 3.615329760455043)
    ,@(5, 9, 18.30282209105937)
    ,@(5, 10, 7.627035134838913)
    ,@(5, 13, 18.90215940323822)
    ,@(5, 15, 19.20765619323754)
    ,@(6, 2, 11.19347403096262)
    ,@(6, 3, 7.861530844577456)
    ,@(6, 5, 21.76597551825019)
    ,@(6, 6, 37.4826293460133)
    ,@(6, 7, 3.615419430619366)
    ,@(6, 9, 18.3098150468513)
    ,@(6, 10, 7.627965435935493)
    ,@(6, 13, 18.89135227871805)
    ,@(6, 15, 19.21322198229094)
    ,@(7, 2, 11.39885823818236)
    ,@(7, 3, 8.013875353994209)
    ,@(7, 5, 21.77328671051905)
    ,@(7, 6, 37.47370729940781)
    ,@(7, 7, 3.614802668948993)
    ,@(7, 9, 18.26172711366804)
    ,@(7, 10, 7.621567127979338)
    ,@(7, 13, 18.966444310341)
    ,@(7, 15, 19.17513361693708)
    ,@(8, 2, 12.26269413550617)
    ,@(8, 3, 8.650736236453902)
    ,@(8, 5, 21.81443806772554)
    ,@(8, 6, 37.45946169094329)
    ,@(8, 7, 3.612220421739792)
    ,@(8, 9, 18.06067145540155)
    ,@(8, 10, 7.594789260069976)
    ,@(8, 13, 19.29928664320772)
    ,@(8, 15, 19.02057502058906)
    ,@(9, 2, 13.79530319051759)
    ,@(9, 3, 9.768687994035789)
    ,@(9, 5, 21.92674001225093)
    ,@(9, 6, 37.52200221372652)
    ,@(9, 7, 3.607654194761944)
    ,@(9, 9, 17.70637116744006)
    ,@(9, 10, 7.547479389903843)
    ,@(9, 13, 19.95283888714672)
    ,@(9, 15, 18.76653526962286)
    ,@(10, 2, 14.81732008779304)
    ,@(10, 3, 10.50800019918054)
    ,@(10, 5, 22.02815206648078)
    ,@(10, 6, 37.62308213386603)
    ,@(10, 7, 3.604599450414921)
    ,@(10, 9, 17.47034765335539)
    ,@(10, 10, 7.515860088937876)
    ,@(10, 13, 20.42959349730621)
    ,@(10, 15, 18.61022343283244)
    ,@(11, 2, 15.25885820249275)
    ,@(11, 3, 10.82621996453286)
    ,@(11, 5, 22.07829931659368)
    ,@(11, 6, 37.6809835411832)
    ,@(11, 7, 3.603274229676375)
    ,@(11, 9, 17.3682305882783)
    ,@(11, 10, 7.502150578307067)
    ,@(11, 13, 20.64493185646186)
    ,@(11, 15, 18.54579189910114)
    ,@(12, 2, 15.42264422198605)
    ,@(12, 3, 10.94410244115853)
    ,@(12, 5, 22.09785698091482)
    ,@(12, 6, 37.70461482467299)
    ,@(12, 7, 3.602781610841402)
    ,@(12, 9, 17.33031578150297)
    ,@(12, 10, 7.497055596405843)
    ,@(12, 13, 20.72618810107729)
    ,@(12, 15, 18.52236102146159)
    ,@(13, 2, 15.38752256646262)
    ,@(13, 3, 10.918831119084)
    ,@(13, 5, 22.09361978164906)
    ,@(13, 6, 37.69944972704211)
    ,@(13, 7, 3.602887296090869)
    ,@(13, 9, 17.33844785272449)
    ,@(13, 10, 7.498148606578344)
    ,@(13, 13, 20.70870194403145)
    ,@(13, 15, 18.52736410292984)
    ,@(14, 2, 15.27240167677814)
    ,@(14, 3, 10.83597082971602)
    ,@(14, 5, 22.07989700968546)
    ,@(14, 6, 37.68289358342628)
    ,@(14, 7, 3.603233517243228)
    ,@(14, 9, 17.365096190681)
    ,@(14, 10, 7.501729479374504)
    ,@(14, 13, 20.65162299504794)
    ,@(14, 15, 18.54384478171595)
    ,@(15, 2, 15.20144072530891)
    ,@(15, 3, 10.78487480259237)
    ,@(15, 5, 22.07156509727364)
    ,@(15, 6, 37.67297424429317)
    ,@(15, 7, 3.603446786181239)
    ,@(15, 9, 17.38151734268565)
    ,@(15, 10, 7.503935422914992)
    ,@(15, 13, 20.61662108995996)
    ,@(15, 15, 18.55406596632681)
    ,@(16, 2, 14.78798926399493)
    ,@(16, 3, 10.48683796038046)
    ,@(16, 5, 22.02495482883153)
    ,@(16, 6, 37.61953730081235)
    ,@(16, 7, 3.604687348544975)
    ,@(16, 9, 17.47712687368285)
    ,@(16, 10, 7.516769567414778)
    ,@(16, 13, 20.41548392953328)
    ,@(16, 15, 18.61456919859809)
    ,@(17, 2, 14.52832024709637)
    ,@(17, 3, 10.29935230977935)
    ,@(17, 5, 21.99738274730674)
    ,@(17, 6, 37.58980290607699)
    ,@(17, 7, 3.605464853766228)
    ,@(17, 9, 17.53712485464174)
    ,@(17, 10, 7.524815276421918)
    ,@(17, 13, 20.29165059846443)
    ,@(17, 15, 18.65340143249264)
    ,@(18, 2, 14.37676694862275)
    ,@(18, 3, 10.18981284669437)
    ,@(18, 5, 21.98190247009149)
    ,@(18, 6, 37.57382307837128)
    ,@(18, 7, 3.605918118406919)
    ,@(18, 9, 17.5721283920115)
    ,@(18, 10, 7.529506447150569)
    ,@(18, 13, 20.22028385872778)
    ,@(18, 15, 18.67636455663924)
    ,@(19, 2, 14.32507779512528)
    ,@(19, 3, 10.15243287032542)
    ,@(19, 5, 21.97672638527946)
    ,@(19, 6, 37.56860562014417)
    ,@(19, 7, 3.606072628991382)
    ,@(19, 9, 17.58406489696635)
    ,@(19, 10, 7.531105715402591)
    ,@(19, 13, 20.19609813069402)
    ,@(19, 15, 18.68424706559113)
    ,@(20, 2, 14.5561904154802)
    ,@(20, 3, 10.31948684572211)
    ,@(20, 5, 22.00027873874493)
    ,@(20, 6, 37.59285205303449)
    ,@(20, 7, 3.605381459725839)
    ,@(20, 9, 17.53068681377838)
    ,@(20, 10, 7.523952229052782)
    ,@(20, 13, 20.30484795659474)
    ,@(20, 15, 18.64920263664881)
    ,@(21, 2, 15.30630851394956)
    ,@(21, 3, 10.86038014194727)
    ,@(21, 5, 22.08391238279598)
    ,@(21, 6, 37.68771032498935)
    ,@(21, 7, 3.603131574017922)
    ,@(21, 9, 17.3572484478657)
    ,@(21, 10, 7.500675074456891)
    ,@(21, 13, 20.66839679231228)
    ,@(21, 15, 18.53897767321523)
    ,@(22, 2, 15.77663274986985)
    ,@(22, 3, 11.1986063040255)
    ,@(22, 5, 22.14187768735516)
    ,@(22, 6, 37.75963961068038)
    ,@(22, 7, 3.60171482547037)
    ,@(22, 9, 17.24829537875308)
    ,@(22, 10, 7.486024451224098)
    ,@(22, 13, 20.90429074573629)
    ,@(22, 15, 18.47258493132148)
    ,@(23, 2, 15.52744916115489)
    ,@(23, 3, 11.01949133813798)
    ,@(23, 5, 22.11064126906272)
    ,@(23, 6, 37.72034414982273)
    ,@(23, 7, 3.602466074328156)
    ,@(23, 9, 17.30604328085225)
    ,@(23, 10, 7.493792456654111)
    ,@(23, 13, 20.77856716010966)
    ,@(23, 15, 18.5075007633862)
    ,@(24, 2, 14.5435973667416)
    ,@(24, 3, 10.31038947608678)
    ,@(24, 5, 21.99896830379236)
    ,@(24, 6, 37.59147006025657)
    ,@(24, 7, 3.605419142647833)
    ,@(24, 9, 17.5335958636621)
    ,@(24, 10, 7.524342208430074)
    ,@(24, 13, 20.29888196571502)
    ,@(24, 15, 18.65109892561653)
    ,@(25, 2, 13.39862276858874)
    ,@(25, 3, 9.480535659224648)
    ,@(25, 5, 21.89301039956489)
    ,@(25, 6, 37.49539064393046)
    ,@(25, 7, 3.608836551850973)
    ,@(25, 9, 17.79794904768518)
    ,@(25, 10, 7.559724402982244)
    ,@(25, 13, 19.77635446864871)
    ,@(25, 15, 18.82996445958117)
)

foreach ($u in $updates) {
    $rowNum = $u[0]
    $colNum = $u[1]
    $value  = $u[2]
    $ws.Cells.Item($rowNum, $colNum).Value = $value
}
